$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove column D entirely (PBASE column dropped)
$ws.Range("D1:D7").EntireColumn.Delete()

# Row 2: corrected raw-value formulas for MeOH only (B) and DMSO only (C)
$ws.Range("B2").Formula = "=-0.05-0.11"
$ws.Range("C2").Formula = "=-0.05-0.1"

# Row 3: corrected raw-value formulas
$ws.Range("B3").Formula = "=0.01-0.14"
$ws.Range("C3").Formula = "=0.02-0.16"

# Row 4: B4 becomes a literal pasted-in value (no formula), with its own font + vertical-centered alignment
$ws.Range("B4").Value = -0.13
$ws.Range("B4").Font.Name = "Calibri"
$ws.Range("B4").Font.Size = 11
$ws.Range("B4").VerticalAlignment = -4108

# Row 5: B5 becomes a literal pasted-in value (no formula), with its own font + 2-decimal number format
$ws.Range("B5").Value = -0.19
$ws.Range("B5").Font.Name = "Calibri"
$ws.Range("B5").Font.Size = 11
$ws.Range("B5").NumberFormat = "0.00"

# Row 6: averages now reference the updated ranges
$ws.Range("B6").Formula = "=AVERAGE(B4:B5)"
$ws.Range("C6").Formula = "=AVERAGE(C2:C3)"

# Row 7: std now references the updated ranges
$ws.Range("B7").Formula = "=STDEV.S(B5:B6)"
$ws.Range("C7").Formula = "=STDEV.S(C2:C3)"

# Selection moved
$ws.Range("F11").Select()

# Page setup: portrait orientation
$ws.PageSetup.Orientation = 1
